$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.148.93"
$ws.Range("E2").Value = "  +2.90%  "
$ws.Range("D3").Value = "2.317.78"
$ws.Range("E3").Value = "  +2.83%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'310.61"
$ws.Range("E5").Value = "  +2.04%  "
$ws.Range("D6").Value = "'101.30"
$ws.Range("E6").Value = "  +6.11%  "
$ws.Range("E7").Value = "  +2.87%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("E9").Value = "  +7.63%  "
$ws.Range("D10").Value = "'36.15"
$ws.Range("E10").Value = "  +4.27%  "
$ws.Range("D11").Value = "'0.0819"
$ws.Range("E11").Value = "  +3.98%  "
$ws.Range("E12").Value = "  +0.79%  "
$ws.Range("D13").Value = "'7.24"
$ws.Range("E13").Value = "  +7.67%  "
$ws.Range("D14").Value = "2.672.28"
$ws.Range("E14").Value = "  +2.55%  "
$ws.Range("D15").Value = "'15.07"
$ws.Range("E15").Value = "  +5.14%  "
$ws.Range("D16").Value = "2.313.07"
$ws.Range("E16").Value = "  +2.34%  "
$ws.Range("D18").Value = "43.092.18"
$ws.Range("E18").Value = "  +3.05%  "
$ws.Range("D19").Value = "'12.63"
$ws.Range("E19").Value = "  +2.64%  "
$ws.Range("D20").Value = "0.0₃0923"
$ws.Range("E20").Value = "  +2.48%  "
$ws.Range("D21").Value = "'6.13"
$ws.Range("E21").Value = "  +3.36%  "
$ws.Range("D22").Value = "'68.63"
$ws.Range("E22").Value = "  +0.81%  "
$ws.Range("D23").Value = "'241.83"
$ws.Range("E23").Value = "  +2.16%  "
$ws.Range("D24").Value = "'2.03"
$ws.Range("E24").Value = "  +5.95%  "
$ws.Range("D25").Value = "'2.63"
$ws.Range("E25").Value = "  +2.96%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("D27").Value = "'24.79"
$ws.Range("E27").Value = "  +5.14%  "
$ws.Range("D28").Value = "'37.53"
$ws.Range("E28").Value = "  +3.62%  "
$ws.Range("D29").Value = "'9.68"
$ws.Range("E29").Value = "  +2.90%  "
$ws.Range("E30").Value = "  -0.23%  "
$ws.Range("D31").Value = "'166.19"
$ws.Range("E31").Value = "  +4.01%  "
$ws.Range("D32").Value = "'5.36"
$ws.Range("E32").Value = "  +3.72%  "
$ws.Range("E33").Value = "  +0.04%  "
$ws.Range("B34").Value = "Celestia"
$ws.Range("C34").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D34").Value = "'17.98"
$ws.Range("E34").Value = "  +6.18%  "
$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").Value = "'3.13"
$ws.Range("E35").Value = "  -1.05%  "
$ws.Range("E36").Value = "  +1.49%  "
$ws.Range("E37").Value = "  +3.88%  "
$ws.Range("E38").Value = "  +0.86%  "
$ws.Range("E39").Value = "  +2.56%  "
$ws.Range("E40").Value = "  +2.47%  "
$ws.Range("D41").Value = "'4.32"
$ws.Range("E41").Value = "  +8.90%  "
$ws.Range("E42").Value = "  +1.24%  "
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").Value = "'19.43"
$ws.Range("E43").Value = "  +4.32%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").Value = "'0.0291"
$ws.Range("E44").Value = "  +3.26%  "
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "1.980.39"
$ws.Range("E45").Value = "  +1.11%  "
$ws.Range("D46").Value = "'3.03"
$ws.Range("E46").Value = "  +4.59%  "
$ws.Range("D47").Value = "'9.85"
$ws.Range("E47").Value = "  +0.10%  "
$ws.Range("D48").Value = "'2.97"
$ws.Range("E48").Value = "  +19.04%  "
$ws.Range("D49").Value = "'55.77"
$ws.Range("E49").Value = "  +6.01%  "
$ws.Range("D50").Value = "2.540.17"
$ws.Range("E50").Value = "  +2.44%  "
$ws.Range("D51").Value = "'1.54"
$ws.Range("E51").Value = "  +3.75%  "
